# worklog on 1/16/2024 sync on clockify
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1/16/2024 block (rows 337-348) ---
# 13:00-15:00 "find jobs"/"Support" -> "Machine learning digit recognition"/"design"
$ws.Range("B343:B345").Value = "Machine learning digit recognition"
$ws.Range("C343:C345").Value = "design"

# 16:00-18:00 "Machine learning digit recognition"/"design" -> "job discussion"/"Support"
$ws.Range("B346:B348").Value = "job discussion"
$ws.Range("C346:C348").Value = "Support"

# --- 1/17/2024 block (rows 349-360) ---
# 17:00-19:00 (row 354) + 13:00-15:00 (355-357) "find jobs"/"Support" -> "Machine learning digit recognition"/"design"
$ws.Range("B354:B357").Value = "Machine learning digit recognition"
$ws.Range("C354:C357").Value = "design"

# 16:00-18:00 (358-360) -> "bussiness close preperation"/"Support"
$ws.Range("B358:B360").Value = "bussiness close preperation"
$ws.Range("C358:C360").Value = "Support"

# --- 1/18/2024 block (rows 361-372) ---
# 17:00-18:00 (366-367) "Machine learning digit recognition"/"design" -> "bussiness close preperation"/"Support"
$ws.Range("B366:B367").Value = "bussiness close preperation"
$ws.Range("C366:C367").Value = "Support"

# Update sheet view scroll position / selection to match saved state
$excel.ActiveWindow.ScrollRow = 318
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C370").Select()
